# Apply the edits described by the diff to the active workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Version: 1.0.1 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# 2. Precondition text (appears in B8, B17, B25, B33 - merged cells, shared string)
$newPrecondition = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B8").Value = $newPrecondition
$ws.Range("B17").Value = $newPrecondition
$ws.Range("B25").Value = $newPrecondition
$ws.Range("B33").Value = $newPrecondition

# 3. TC1 expected result - fix accents (numero->número, diaria->diária, diarias->diárias)
$ws.Range("D11").Value = "SYSTEM Exibe a lista de diárias (solicitações) aptas para pagamento ordenado pelo número da diária em ordem crescente. Exibe esta lista de diárias também ordenada pela data de chegada da solicitação na fase de liquidação (após registrar o empenho)."

# 4. TC2 expected result - add trailing period
$ws.Range("D20").Value = "SYSTEM Apresenta a tela de Registrar Liquidações."

# 5. TC3 expected result - remove redundant "o nome"
$ws.Range("D28").Value = "SYSTEM Atualiza a lista de registros de solicitações, onde deverá constar o nome do usuário logado (que se atribuiu como responsável pela liquidação) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

# 6. TC4 expected result - add trailing period
$ws.Range("D36").Value = "SYSTEM Apresenta a tela de Detalhar Diárias."
